$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 137: George Lucas: A Life ---
$ws.Range("A136:I136").Copy()
$ws.Range("A137:I137").PasteSpecial(-4122)
$ws.Cells.Item(137, 1).Value = "George Lucas: A Life"
$ws.Cells.Item(137, 2).Value = "Brian Jay Jones"
$ws.Cells.Item(137, 3).Value = 44130
$ws.Cells.Item(137, 4).Value = 44132
$ws.Cells.Item(137, 5).Value = "biography;history;george lucas;star wars;hollywood;cinema;indiana jones"
$ws.Cells.Item(137, 6).Value = "Audio"
$ws.Cells.Item(137, 7).Value = "18 Hours 18 Mins"
$ws.Cells.Item(137, 8).Value = 4
$ws.Cells.Item(137, 9).Value = $true

# --- Row 138: Buffett ---
$ws.Range("A137:I137").Copy()
$ws.Range("A138:I138").PasteSpecial(-4122)
$ws.Cells.Item(138, 1).Value = "Buffett"
$ws.Cells.Item(138, 2).Value = "Roger Lowenstein"
$ws.Cells.Item(138, 3).Value = 44100
$ws.Cells.Item(138, 4).Value = 44133
$ws.Cells.Item(138, 5).Value = "biography;warren buffett;finance;investing"
$ws.Cells.Item(138, 6).Value = "Hard Copy"
$ws.Cells.Item(138, 7).Value = "423 Pages"
$ws.Cells.Item(138, 8).Value = 3
$ws.Cells.Item(138, 9).Value = $true

# Move selection to the next empty row, like a user would after data entry
$ws.Range("A139").Select() | Out-Null
